$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.560.84'
$ws.Range("E2").Value = '  -1.64%  '

$ws.Range("D3").Value = '2.619.94'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.12'
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.02'
$ws.Range("E6").Value = '  +0.78%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  +1.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.86'
$ws.Range("E9").Value = '  +5.05%  '

$ws.Range("E10").Value = '  -2.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.333'
$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("E12").Value = '  +1.23%  '

$ws.Range("D13").Value = '3.085.78'
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("D14").Value = '58.481.38'
$ws.Range("E14").Value = '  -1.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.70'
$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("D16").Value = '2.645.23'
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("E17").Value = '  -1.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.42'
$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '334.42'
$ws.Range("E19").Value = '  -2.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.15'
$ws.Range("E20").Value = '  +0.50%  '

$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.26'
$ws.Range("E23").Value = '  -1.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.419'
$ws.Range("E24").Value = '  +2.45%  '

$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.62%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  -1.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.10'
$ws.Range("E27").Value = '  -1.77%  '

$ws.Range("E28").Value = '  -1.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  -2.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.89'
$ws.Range("E31").Value = '  +0.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.87'
$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '151.67'
$ws.Range("E33").Value = '  +1.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.90'
$ws.Range("E34").Value = '  -2.07%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").Value = '  -2.07%  '

$ws.Range("B36").Value = 'SuiNetwork'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.847'
$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.815'
$ws.Range("E37").Value = '  -1.34%  '

$ws.Range("E38").Value = '  -2.92%  '

$ws.Range("E39").Value = '  +0.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '281.68'
$ws.Range("E40").Value = '  +3.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.594'

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0940'
$ws.Range("E44").Value = '  -1.16%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.99'
$ws.Range("E45").Value = '  +2.07%  '

$ws.Range("E46").Value = '  +0.84%  '

$ws.Range("D48").Value = '1.946.25'
$ws.Range("E48").Value = '  +0.33%  '

$ws.Range("E49").Value = '  -1.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.88'
$ws.Range("E50").Value = '  -3.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '114.10'
$ws.Range("E51").Value = '  +1.96%  '
